# edit.ps1 - applies the small wording/layout/table-style tweaks described
# in the commit "small changes in the presentation".
#
# Helper: replace an exact substring inside a shape's TextRange while
# preserving run boundaries/formatting (by locating the substring via the
# COM-side TextRange.Text/IndexOf and rewriting exactly that character
# span with TextRange.Characters(start, length).Text = ...).
function Set-ExactTextRange($TextRange, $OldText, $NewText) {
    $full = $TextRange.Text
    $idx0 = $full.IndexOf($OldText)
    if ($idx0 -lt 0) {
        throw "Text not found: $OldText"
    }
    $span = $TextRange.Characters($idx0 + 1, $OldText.Length)
    $span.Text = $NewText
}

$p = $ppt.ActivePresentation

# --- Slide 2 ("Modelling"): wording + roundRect reposition ---------------
$slide2 = $p.Slides.Item(2)

# Shape 70 ("Aircraft (AC): Main module of the analysis composed by: ...").
# This shape auto-sizes (<a:spAutoFit/>), and this runtime's text-measurement
# recomputes the box height as a side effect of any TextRange write (even when
# the replacement text is the same length) - the source diff doesn't touch the
# shape's size, so explicitly restore its original authored height afterwards
# (196.346459 pt round-trips to the original 2493600 EMU).
$shpText = $slide2.Shapes.Item(4)
Set-ExactTextRange $shpText.TextFrame.TextRange ": Main module of the analysis composed by:" ": Main module of the analysis composed of:"
$shpText.Height = 196.346459

# Shape 71 (round-rect outline behind the bullet text) moves slightly.
$shpRect = $slide2.Shapes.Item(5)
$shpRect.Left = 4475600 / 12700.0
$shpRect.Top = 1266600 / 12700.0

# --- Slide 4: table style id -----------------------------------------------
$slide4 = $p.Slides.Item(4)
$tblShape = $slide4.Shapes.Item(7)
$tblShape.Table.ApplyStyle("{FF7EBE3A-FCFE-4183-ABAB-3D2FE987B653}")

# --- Slide 6: wording tweaks -----------------------------------------------
$slide6 = $p.Slides.Item(6)
$shpPara = $slide6.Shapes.Item(3)
$tr6 = $shpPara.TextFrame.TextRange

Set-ExactTextRange $tr6 "Below are reported the " "In the graphs are reported the "

Set-ExactTextRange $tr6 " behavior related to the duration of the simulation in monitored mode, with exponential time generation" " behaviors related to the duration of the simulation in monitored mode, with exponential time generation"
